$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183 (existing rows 183-222 shift down to 184-223)
$ws.Rows.Item(183).Insert()

# Populate the new row 183 with the weekly record
$ws.Range("A183").Value = 5
$ws.Range("B183").Value = "Macroferia Regional de Talca"
$ws.Range("C183").Value = "Maule"
$ws.Range("D183").Value = 44543
$ws.Range("E183").Value = 7
$ws.Range("F183").Value = 100114014
$ws.Range("G183").Value = "Betarraga"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 5000
$ws.Range("K183").Value = 500
$ws.Range("L183").Value = 500
$ws.Range("M183").Value = 500
$ws.Range("N183").Value = "$/paquete 5 unidades"
$ws.Range("O183").Value = "Región del Maule"
$ws.Range("P183").Value = 100
$ws.Range("Q183").Value = 5
$ws.Range("R183").Value = "Hortaliza"
